$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2024-04-22 Monday" "2024-04-23 Tuesday"

Replace-Text "51×48=" "54×24="
Replace-Text "49×12=" "90×84="
Replace-Text "31×44=" "97×22="
Replace-Text "55×46=" "99×41="
Replace-Text "78×43=" "62×75="

Replace-Text "22×16=" "49×91="
Replace-Text "89×32=" "93×20="
Replace-Text "55×16=" "23×75="
Replace-Text "19×61=" "95×45="
Replace-Text "96×13=" "56×33="

Replace-Text "61×11=" "50×89="
Replace-Text "93×40=" "25×47="
Replace-Text "97×44=" "62×21="
Replace-Text "56×17=" "89×89="
Replace-Text "34×29=" "95×99="

Replace-Text "83×80=" "31×51="
Replace-Text "84×54=" "15×66="
Replace-Text "93×63=" "25×21="
Replace-Text "45×96=" "50×33="
Replace-Text "45×67=" "31×47="

Replace-Text "26×68=" "45×41="
Replace-Text "85×44=" "75×55="
Replace-Text "37×63=" "31×36="
Replace-Text "72×46=" "20×86="
Replace-Text "88×86=" "24×23="
